$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.267.16'
$ws.Cells.Item(2, 5).Value = '  -3.83%  '
$ws.Cells.Item(3, 4).Value = '2.246.93'
$ws.Cells.Item(3, 5).Value = '  -4.51%  '
$ws.Cells.Item(4, 5).Value = '  -0.09%  '
$c = $ws.Cells.Item(5, 4)
$c.Value = "'232.77"
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -3.23%  '
$ws.Cells.Item(6, 5).Value = '  -5.45%  '
$c = $ws.Cells.Item(7, 4)
$c.Value = "'70.52"
$c.Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -3.92%  '
$ws.Cells.Item(8, 5).Value = '  +0.05%  '
$c = $ws.Cells.Item(9, 4)
$c.Value = "'0.566"
$c.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -6.02%  '
$ws.Cells.Item(10, 5).Value = '  -0.85%  '
$c = $ws.Cells.Item(11, 4)
$c.Value = "'58.59"
$c.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +0.34%  '
$c = $ws.Cells.Item(12, 4)
$c.Value = "'36.38"
$c.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +7.73%  '
$ws.Cells.Item(13, 5).Value = '  -2.85%  '
$ws.Cells.Item(14, 5).Value = '  -6.01%  '
$ws.Cells.Item(15, 4).Value = '2.575.89'
$ws.Cells.Item(15, 5).Value = '  -4.77%  '
$c = $ws.Cells.Item(16, 4)
$c.Value = "'15.09"
$c.Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -8.21%  '
$ws.Cells.Item(17, 5).Value = '  -4.70%  '
$ws.Cells.Item(18, 4).Value = '2.245.32'
$ws.Cells.Item(18, 5).Value = '  -4.50%  '
$ws.Cells.Item(19, 4).Value = '42.108.46'
$ws.Cells.Item(19, 5).Value = '  -3.92%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0989'
$ws.Cells.Item(20, 5).Value = '  -3.61%  '
$c = $ws.Cells.Item(21, 4)
$c.Value = "'6.27"
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -6.61%  '
$c = $ws.Cells.Item(22, 4)
$c.Value = "'73.62"
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -5.35%  '
$c = $ws.Cells.Item(23, 4)
$c.Value = "'238.04"
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -7.47%  '
$c = $ws.Cells.Item(24, 4)
$c.Value = "'2.02"
$c.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +4.55%  '
$c = $ws.Cells.Item(25, 4)
$c.Value = "'1.00"
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -0.23%  '
$ws.Cells.Item(26, 5).Value = '  -2.57%  '
$c = $ws.Cells.Item(27, 4)
$c.Value = "'2.36"
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -5.78%  '
$c = $ws.Cells.Item(28, 4)
$c.Value = "'10.15"
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -4.92%  '
$c = $ws.Cells.Item(29, 4)
$c.Value = "'2.16"
$c.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -9.19%  '
$c = $ws.Cells.Item(30, 4)
$c.Value = "'168.40"
$c.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -5.19%  '
$c = $ws.Cells.Item(31, 4)
$c.Value = "'20.75"
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -8.51%  '
$ws.Cells.Item(32, 5).Value = '  -7.31%  '
$ws.Cells.Item(33, 5).Value = '  -6.51%  '
$ws.Cells.Item(34, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Cells.Item(34, 4)
$c.Value = "'5.41"
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -1.27%  '
$ws.Cells.Item(35, 2).Value = 'Hedera'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(35, 4)
$c.Value = "'0.0722"
$c.Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -4.63%  '
$ws.Cells.Item(36, 5).Value = '  -7.56%  '
$ws.Cells.Item(37, 5).Value = '  -5.06%  '
$c = $ws.Cells.Item(38, 4)
$c.Value = "'22.34"
$c.Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +17.20%  '
$c = $ws.Cells.Item(39, 4)
$c.Value = "'6.11"
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -5.16%  '
$ws.Cells.Item(40, 5).Value = '  -5.42%  '
$c = $ws.Cells.Item(41, 4)
$c.Value = "'0.0269"
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -3.33%  '
$c = $ws.Cells.Item(42, 4)
$c.Value = "'67.47"
$c.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +0.20%  '
$c = $ws.Cells.Item(43, 4)
$c.Value = "'5.02"
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -4.32%  '
$c = $ws.Cells.Item(44, 4)
$c.Value = "'9.07"
$c.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -1.21%  '
$ws.Cells.Item(45, 5).Value = '  -9.34%  '
$ws.Cells.Item(46, 5).Value = '  -6.66%  '
$c = $ws.Cells.Item(47, 4)
$c.Value = "'1.00"
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -0.02%  '
$ws.Cells.Item(48, 2).Value = 'SynthetixNetwork'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$c = $ws.Cells.Item(48, 4)
$c.Value = "'4.43"
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +7.48%  '
$ws.Cells.Item(49, 2).Value = 'Celestia'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$c = $ws.Cells.Item(49, 4)
$c.Value = "'10.28"
$c.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +8.06%  '
$ws.Cells.Item(50, 2).Value = 'NEARProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Cells.Item(50, 4)
$c.Value = "'2.38"
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -5.53%  '
$ws.Cells.Item(51, 5).Value = '  -5.85%  '
